$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (DAMSLTag, DialogAct)
$updates = @{
    2   = @("aa", "Agree/Accept")
    4   = @("qy", "Yes-No-Question")
    10  = @("sd", "Statement-non-opinion")
    17  = @("sd", "Statement-non-opinion")
    18  = @("sd", "Statement-non-opinion")
    22  = @("sv", "Statement-opinion")
    23  = @("qy", "Yes-No-Question")
    59  = @("ba", "Appreciation")
    61  = @("sd", "Statement-non-opinion")
    63  = @("aa", "Agree/Accept")
    73  = @("sv", "Statement-opinion")
    74  = @("sv", "Statement-opinion")
    77  = @("sd", "Statement-non-opinion")
    111 = @("aa", "Agree/Accept")
    133 = @("b", "Acknowledge (Backchannel)")
    137 = @("%", "Uninterpretable")
    161 = @("sd", "Statement-non-opinion")
    172 = @("sd", "Statement-non-opinion")
    178 = @("sd", "Statement-non-opinion")
    186 = @("%", "Uninterpretable")
    192 = @("b", "Acknowledge (Backchannel)")
    208 = @("ba", "Appreciation")
    209 = @("sv", "Statement-opinion")
    210 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
